# Generate Report for Handoff
#
# Updates the localization-status report to reflect a fresh handoff:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - The handoff timestamps are refreshed
#   - The (now shorter) status column is narrowed to fit its new content
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Status column on all three sheets ---------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Refreshed handoff timestamps ---------------------------------------
$overview.Range("G2").Value = "2016-08-26 16:59:56"
$dede.Range("H2").Value = "2016-08-26 16:59:56"
$zhcn.Range("H2").Value = "2016-08-26 16:59:52"

# --- Narrow the status columns to fit the shorter text -------------------
# (engine rounds ColumnWidth to an internal 1/6-character pixel grid, so we
#  feed the value whose rounded result lands closest to the target width)
$narrowWidth = 16.333333333333332
$overview.Range("E1").ColumnWidth = $narrowWidth
$overview.Range("F1").ColumnWidth = $narrowWidth
$zhcn.Range("C1").ColumnWidth = $narrowWidth
$dede.Range("C1").ColumnWidth = $narrowWidth

Write-Host "Handoff report regenerated"
